# The deck originally has the "Integral" theme driving the slide master
# (ppt/theme/theme1.xml) while the "Office Theme" colours sit unused in
# ppt/theme/theme2.xml (the notes master's theme). The authored change
# swaps the two themes' content so the presentation's visible design
# becomes the stock "Office Theme" colour palette.
#
# The font scheme and format scheme (fills/lines/effects) are byte-for-byte
# identical between the two themes, so only the 12 theme colours need to
# change. We drive that through the slide master's theme colour scheme,
# which every slide shares (there is only one design/master in this deck).

$p = $ppt.ActivePresentation

# Office Theme's 12 standard theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) expressed as COM RGB (0xBBGGRR) longs.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeThemeColors[$i - 1]
}
